$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("WillKoehrsen", "Data-Analysis", 42809.80076388889, 44231.51909722222, 0, 376839, "Jupyter Notebook", 6, 4, 350, 4051, 3285, 56, 367, 15, 3, 0, 0, "Data Science Using Python"),
    @("microsoft", "tensorwatch", 43600.35386574074, 44299.40557870371, 0, 22242, "Jupyter Notebook", 3, 10000, 104, 3182, 346, 65, 112, 14, 1, 7, 0, "Debugging, monitoring and visualization for Python Machine Learning and Data Science"),
    @("nborwankar", "LearnDataScience", 41481.10930555555, 43760.41459490741, 0, 148439, "Jupyter Notebook", 3, 7, 371, 2673, 1647, 1, 64, 19, 2, 0, 0, "Open Content for self-directed learning in data science"),
    @("ml-tooling", "ml-workspace", 43612.70503472222, 44480.09324074074, 0, 13023, "Jupyter Notebook", 7, 9, 62, 2287, 314, 73, 845, 29, 4, 34, 4, "🛠 All-in-one web-based IDE specialized for machine learning and data science."),
    @("ShuaiW", "data-science-question-answer", 43115.06939814815, 44105.20944444444, 0, 2390, "Jupyter Notebook", 1, 2, 88, 1980, 603, 1, 17, 4, 1, 0, 0, "A repo for data science related questions and answers")
)

$startRow = 7
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $cell = $ws.Cells.Item($r, $c + 1)
        $cell.Value = $row[$c]
    }
    # Date-formatted columns C (3) and D (4) mirror the style used in rows above
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
